# Update "想去人数" (interested count, column F) for a handful of events
# that appear on both the "展览" and "全部类型" sheets.
#
# Sheet "展览":     F8 1732->1733, F16 1553->1554, F17 7177->7178,
#                    F19 7318->7319, F21 12->13, F22 5629->5631, F43 1117->1118
# Sheet "全部类型": F6 1732->1733, F14 1553->1554, F18 7177->7178,
#                    F20 7318->7319, F22 12->13, F23 5629->5631, F46 1117->1118

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F8").Value = 1733
$ws1.Range("F16").Value = 1554
$ws1.Range("F17").Value = 7178
$ws1.Range("F19").Value = 7319
$ws1.Range("F21").Value = 13
$ws1.Range("F22").Value = 5631
$ws1.Range("F43").Value = 1118

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F6").Value = 1733
$ws4.Range("F14").Value = 1554
$ws4.Range("F18").Value = 7178
$ws4.Range("F20").Value = 7319
$ws4.Range("F22").Value = 13
$ws4.Range("F23").Value = 5631
$ws4.Range("F46").Value = 1118
